$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.4050618212753369
$ws.Range("J4").Value = 0.4928632049788898
$ws.Range("K4").Value = 0.4256430241151026
$ws.Range("L4").Value = 2.652482080263422
